# Add BAU CCS subsidy
#
# Insert a new "GRA-CCSsubsidy" worksheet right after "GRA-fuelsubsidy" (and
# before "GRA-ntnldebtinterest"), by duplicating the GRA-fuelsubsidy sheet
# (same layout / formulas / values) and renaming the copy.

$wb = $excel.ActiveWorkbook

$srcSheet = $wb.Worksheets.Item("GRA-fuelsubsidy")

# Copy the sheet, placing the new copy immediately after the source sheet.
$srcSheet.Copy($null, $srcSheet)

# Excel names the freshly-copied sheet "GRA-fuelsubsidy (2)"; rename it.
$newSheet = $wb.Worksheets.Item("GRA-fuelsubsidy (2)")
$newSheet.Name = "GRA-CCSsubsidy"
